# Update leve/profit calculation columns (currentAveragePrice, NQ/HQ prices,
# leve buy prices, and NQ/HQ profit) across all eight crafting-class sheets.
# Values refresh from an external market-data snapshot (scheduled runner).

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 558.7143
$ws.Range("I111").Value = 343.33334
$ws.Range("J111").Value = 720.25
$ws.Range("K111").Value = 1030.00002
$ws.Range("L111").Value = 2160.75
$ws.Range("M111").Value = 2036.99998
$ws.Range("N111").Value = -8294.75
$ws.Range("H132").Value = 15076227
$ws.Range("I132").Value = 17858296
$ws.Range("J132").Value = 912967.4399999999
$ws.Range("K132").Value = 53574888
$ws.Range("L132").Value = 2738902.32
$ws.Range("M132").Value = -53572358
$ws.Range("N132").Value = -2743962.32

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 837.1667
$ws.Range("I2").Value = 952.5
$ws.Range("K2").Value = 952.5
$ws.Range("M2").Value = -839.5
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 421.57144
$ws.Range("I5").Value = 390.2
$ws.Range("K5").Value = 390.2
$ws.Range("M5").Value = -278.2
$ws.Range("H9").Value = 18251
$ws.Range("J9").Value = 18251
$ws.Range("L9").Value = 18251
$ws.Range("N9").Value = -18591
$ws.Range("H20").Value = 18251
$ws.Range("J20").Value = 18251
$ws.Range("L20").Value = 18251
$ws.Range("N20").Value = -18791
$ws.Range("H23").Value = 17514.4
$ws.Range("J23").Value = 17514.4
$ws.Range("L23").Value = 17514.4
$ws.Range("N23").Value = -18032.4
$ws.Range("H24").Value = 36871
$ws.Range("J24").Value = 36871
$ws.Range("L24").Value = 36871
$ws.Range("N24").Value = -37619
$ws.Range("H37").Value = 24637.666
$ws.Range("I37").Value = 15417
$ws.Range("J37").Value = 27272.143
$ws.Range("K37").Value = 15417
$ws.Range("L37").Value = 27272.143
$ws.Range("M37").Value = -15144
$ws.Range("N37").Value = -27818.143
$ws.Range("H44").Value = 41152
$ws.Range("J44").Value = 41152
$ws.Range("L44").Value = 41152
$ws.Range("N44").Value = -42128
$ws.Range("H55").Value = 43364
$ws.Range("J55").Value = 43364
$ws.Range("L55").Value = 43364
$ws.Range("N55").Value = -43994
$ws.Range("H63").Value = 13854401
$ws.Range("I63").Value = 17316502
$ws.Range("K63").Value = 17316502
$ws.Range("M63").Value = -17315816
$ws.Range("H66").Value = 13854401
$ws.Range("I66").Value = 17316502
$ws.Range("K66").Value = 86582510
$ws.Range("M66").Value = -86579078
$ws.Range("H74").Value = 2694.04
$ws.Range("I74").Value = 2511.2273
$ws.Range("K74").Value = 2511.2273
$ws.Range("M74").Value = -1637.2273
$ws.Range("H77").Value = 2694.04
$ws.Range("I77").Value = 2511.2273
$ws.Range("K77").Value = 12556.1365
$ws.Range("M77").Value = -8188.136500000001
$ws.Range("H100").Value = 36871
$ws.Range("J100").Value = 36871
$ws.Range("L100").Value = 36871
$ws.Range("N100").Value = -39035
$ws.Range("H116").Value = 837.1667
$ws.Range("I116").Value = 952.5
$ws.Range("K116").Value = 952.5
$ws.Range("M116").Value = 1341.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 837.1667
$ws.Range("I3").Value = 952.5
$ws.Range("K3").Value = 952.5
$ws.Range("M3").Value = -838.5
$ws.Range("H4").Value = 421.57144
$ws.Range("I4").Value = 390.2
$ws.Range("K4").Value = 390.2
$ws.Range("M4").Value = -275.2
$ws.Range("H15").Value = 31000
$ws.Range("J15").Value = 31000
$ws.Range("L15").Value = 31000
$ws.Range("N15").Value = -31454
$ws.Range("H19").Value = 35005
$ws.Range("J19").Value = 35005
$ws.Range("L19").Value = 35005
$ws.Range("N19").Value = -35351
$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 50000
$ws.Range("L35").Value = 50000
$ws.Range("N35").Value = -50620
$ws.Range("H82").Value = 22782
$ws.Range("I82").Value = 3742.5
$ws.Range("J82").Value = 29705.455
$ws.Range("K82").Value = 3742.5
$ws.Range("L82").Value = 29705.455
$ws.Range("M82").Value = -3359.5
$ws.Range("N82").Value = -30471.455
$ws.Range("H85").Value = 22782
$ws.Range("I85").Value = 3742.5
$ws.Range("J85").Value = 29705.455
$ws.Range("K85").Value = 3742.5
$ws.Range("L85").Value = 29705.455
$ws.Range("M85").Value = -2416.5
$ws.Range("N85").Value = -32357.455
$ws.Range("H99").Value = 4539.8
$ws.Range("J99").Value = 4539.8
$ws.Range("L99").Value = 4539.8
$ws.Range("N99").Value = -7535.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5548.7856
$ws.Range("I17").Value = 464.8889
$ws.Range("J17").Value = 14699.8
$ws.Range("K17").Value = 464.8889
$ws.Range("L17").Value = 14699.8
$ws.Range("M17").Value = -290.8889
$ws.Range("N17").Value = -15047.8
$ws.Range("H25").Value = 9940.75
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 9940.75
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 9940.75
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -10288.75
$ws.Range("H41").Value = 35601.715
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 35601.715
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 35601.715
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -36457.715
$ws.Range("H50").Value = 29460.666
$ws.Range("J50").Value = 29460.666
$ws.Range("L50").Value = 29460.666
$ws.Range("N50").Value = -30710.666
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290
$ws.Range("H60").Value = 31632.182
$ws.Range("J60").Value = 31632.182
$ws.Range("L60").Value = 31632.182
$ws.Range("N60").Value = -32654.182
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696
$ws.Range("H74").Value = 36174.89
$ws.Range("J74").Value = 36174.89
$ws.Range("L74").Value = 36174.89
$ws.Range("N74").Value = -37922.89
$ws.Range("H77").Value = 36174.89
$ws.Range("J77").Value = 36174.89
$ws.Range("L77").Value = 108524.67
$ws.Range("N77").Value = -117260.67

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1216.659
$ws.Range("I5").Value = 286.86365
$ws.Range("K5").Value = 860.59095
$ws.Range("M5").Value = -748.59095
$ws.Range("H135").Value = 1216.659
$ws.Range("I135").Value = 286.86365
$ws.Range("K135").Value = 2581.77285
$ws.Range("M135").Value = -46.77285000000029

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 14925
$ws.Range("J18").Value = 14925
$ws.Range("L18").Value = 14925
$ws.Range("N18").Value = -15511
$ws.Range("H43").Value = 16874.666
$ws.Range("J43").Value = 26886.727
$ws.Range("L43").Value = 26886.727
$ws.Range("N43").Value = -27188.727
$ws.Range("H58").Value = 15883.333
$ws.Range("I58").Value = 9575
$ws.Range("J58").Value = 28500
$ws.Range("K58").Value = 9575
$ws.Range("L58").Value = 28500
$ws.Range("M58").Value = -9298
$ws.Range("N58").Value = -29054
$ws.Range("H64").Value = 35317.145
$ws.Range("J64").Value = 35317.145
$ws.Range("L64").Value = 35317.145
$ws.Range("N64").Value = -35813.145
$ws.Range("H67").Value = 35317.145
$ws.Range("J67").Value = 35317.145
$ws.Range("L67").Value = 35317.145
$ws.Range("N67").Value = -37033.145
$ws.Range("H80").Value = 25002282
$ws.Range("I80").Value = 35716204
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 35716204
$ws.Range("L80").Value = 3133.3333
$ws.Range("M80").Value = -35715206
$ws.Range("N80").Value = -5129.3333
$ws.Range("H83").Value = 25002282
$ws.Range("I83").Value = 35716204
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 178581020
$ws.Range("L83").Value = 15666.6665
$ws.Range("M83").Value = -178576028
$ws.Range("N83").Value = -25650.6665
$ws.Range("H132").Value = 2287.4634
$ws.Range("I132").Value = 1458.6875
$ws.Range("J132").Value = 5234.222
$ws.Range("K132").Value = 4376.0625
$ws.Range("L132").Value = 15702.666
$ws.Range("M132").Value = -1846.0625
$ws.Range("N132").Value = -20762.666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5730.7646
$ws.Range("I122").Value = 2935.3333
$ws.Range("K122").Value = 8805.999899999999
$ws.Range("M122").Value = -6355.999899999999
$ws.Range("H132").Value = 11194.892
$ws.Range("I132").Value = 11917.167
$ws.Range("J132").Value = 9861.462
$ws.Range("K132").Value = 35751.501
$ws.Range("L132").Value = 29584.386
$ws.Range("M132").Value = -33221.501
$ws.Range("N132").Value = -34644.386

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8334833.5
$ws.Range("I132").Value = 602.2414
$ws.Range("J132").Value = 30306898
$ws.Range("K132").Value = 1806.7242
$ws.Range("L132").Value = 90920694
$ws.Range("M132").Value = 723.2757999999999
$ws.Range("N132").Value = -90925754
